$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets("ALC")
$ws.Range("H13").Value = 500
$ws.Range("J13").Value = 500
$ws.Range("L13").Value = 500
$ws.Range("N13").Value = -838

$ws.Range("H55").Value = 190.125
$ws.Range("I55").Value = 141.66667
$ws.Range("J55").Value = 219.2
$ws.Range("K55").Value = 141.66667
$ws.Range("L55").Value = 219.2
$ws.Range("M55").Value = 72.33332999999999
$ws.Range("N55").Value = -647.2

$ws.Range("H107").Value = 986.3333
$ws.Range("I107").Value = 617
$ws.Range("J107").Value = 2094.3333
$ws.Range("K107").Value = 617
$ws.Range("L107").Value = 2094.3333
$ws.Range("M107").Value = 1303
$ws.Range("N107").Value = -5934.3333

$ws.Range("H117").Value = 38200
$ws.Range("J117").Value = 38200
$ws.Range("L117").Value = 38200
$ws.Range("N117").Value = -47378

$ws.Range("H132").Value = 38000
$ws.Range("I132").Value = 38000
$ws.Range("K132").Value = 114000
$ws.Range("M132").Value = -111470

$ws = $wb.Sheets("ARM")
$ws.Range("H3").Value = 104.75
$ws.Range("I3").Value = 116.333336
$ws.Range("J3").Value = 70
$ws.Range("K3").Value = 116.333336
$ws.Range("L3").Value = 70
$ws.Range("M3").Value = -1.333336000000003
$ws.Range("N3").Value = -300

$ws.Range("H32").Value = 438.93
$ws.Range("I32").Value = 465
$ws.Range("J32").Value = 247.75
$ws.Range("K32").Value = 465
$ws.Range("L32").Value = 247.75
$ws.Range("M32").Value = -178
$ws.Range("N32").Value = -821.75

$ws.Range("H45").Value = 1312.375
$ws.Range("I45").Value = 1312.375
$ws.Range("K45").Value = 1312.375
$ws.Range("M45").Value = -935.375

$ws = $wb.Sheets("BSM")
$ws.Range("H5").Value = 15246
$ws.Range("I5").Value = 587
$ws.Range("J5").Value = 54336.668
$ws.Range("K5").Value = 587
$ws.Range("L5").Value = 54336.668
$ws.Range("M5").Value = -474
$ws.Range("N5").Value = -54562.668

$ws.Range("H24").Value = 15000.5
$ws.Range("I24").Value = 15000.5
$ws.Range("K24").Value = 15000.5
$ws.Range("M24").Value = -14765.5

$ws.Range("H36").Value = 6685
$ws.Range("I36").Value = 4783.5713
$ws.Range("K36").Value = 4783.5713
$ws.Range("M36").Value = -4249.5713

$ws.Range("H107").Value = 2413.72
$ws.Range("I107").Value = 1922.8572
$ws.Range("J107").Value = 4990.75
$ws.Range("K107").Value = 1922.8572
$ws.Range("L107").Value = 4990.75
$ws.Range("M107").Value = -2.857199999999921
$ws.Range("N107").Value = -8830.75

$ws.Range("H141").Value = 43232.332
$ws.Range("I141").Value = 19854.5
$ws.Range("K141").Value = 19854.5
$ws.Range("M141").Value = -14674.5

$ws = $wb.Sheets("CRP")
$ws.Range("H31").Value = 7360.8184
$ws.Range("J31").Value = 7607.3125
$ws.Range("L31").Value = 7607.3125
$ws.Range("N31").Value = -8197.3125

$ws.Range("H34").Value = 7360.8184
$ws.Range("J34").Value = 7607.3125
$ws.Range("L34").Value = 7607.3125
$ws.Range("N34").Value = -8011.3125

$ws.Range("H39").Value = 3499.5
$ws.Range("I39").Value = 3499.5
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 3499.5
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -3108.5
$ws.Range("N39").ClearContents()

$ws.Range("H49").Value = 3499.5
$ws.Range("I49").Value = 3499.5
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 3499.5
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -3317.5
$ws.Range("N49").ClearContents()

$ws.Range("H58").Value = 3423.9688
$ws.Range("I58").Value = 3185.238
$ws.Range("K58").Value = 3185.238
$ws.Range("M58").Value = -2982.238

$ws.Range("H95").Value = 9974.666999999999
$ws.Range("J95").Value = 9974.666999999999
$ws.Range("L95").Value = 9974.666999999999
$ws.Range("N95").Value = -15466.667

$ws.Range("H99").Value = 5072.222
$ws.Range("I99").Value = 5206.25
$ws.Range("K99").Value = 5206.25
$ws.Range("M99").Value = -3708.25

$ws.Range("H117").Value = 23385.8
$ws.Range("J117").Value = 34997
$ws.Range("L117").Value = 34997
$ws.Range("N117").Value = -44175

$ws.Range("H126").Value = 5072.222
$ws.Range("I126").Value = 5206.25
$ws.Range("K126").Value = 15618.75
$ws.Range("M126").Value = -13148.75

$ws.Range("H131").Value = 29846.375
$ws.Range("J131").Value = 32395.857
$ws.Range("L131").Value = 32395.857
$ws.Range("N131").Value = -42475.857

$ws.Range("H132").Value = 9996.333000000001
$ws.Range("I132").Value = 9996.333000000001
$ws.Range("K132").Value = 29988.999
$ws.Range("M132").Value = -27458.999

$ws.Range("H134").Value = 5811.2856
$ws.Range("I134").Value = 6005.5
$ws.Range("K134").Value = 18016.5
$ws.Range("M134").Value = -15481.5

$ws.Range("H136").Value = 3423.9688
$ws.Range("I136").Value = 3185.238
$ws.Range("K136").Value = 9555.714
$ws.Range("M136").Value = -7005.714

$ws.Range("H141").Value = 102133
$ws.Range("J141").Value = 130904.29
$ws.Range("L141").Value = 130904.29
$ws.Range("N141").Value = -141264.29

$ws = $wb.Sheets("CUL")
$ws.Range("H37").Value = 102499.836
$ws.Range("J37").Value = 102499.836
$ws.Range("L37").Value = 307499.508
$ws.Range("N37").Value = -307723.508

$ws.Range("H70").Value = 436.66666
$ws.Range("J70").Value = 399
$ws.Range("L70").Value = 1197
$ws.Range("N70").Value = -1827

$ws.Range("H73").Value = 436.66666
$ws.Range("J73").Value = 399
$ws.Range("L73").Value = 1197
$ws.Range("N73").Value = -3381

$ws.Range("H128").Value = 139995
$ws.Range("I128").Value = 139995
$ws.Range("K128").Value = 419985
$ws.Range("M128").Value = -415005

$ws = $wb.Sheets("GSM")
$ws.Range("H17").Value = 1541.6
$ws.Range("J17").Value = 2500
$ws.Range("L17").Value = 2500
$ws.Range("N17").Value = -2836

$ws.Range("H123").Value = 52500
$ws.Range("I123").Value = 53333.332
$ws.Range("J123").Value = 51666.668
$ws.Range("K123").Value = 53333.332
$ws.Range("L123").Value = 51666.668
$ws.Range("M123").Value = -50883.332
$ws.Range("N123").Value = -56566.668

$ws.Range("H126").Value = 49761.855
$ws.Range("I126").Value = 60752.06
$ws.Range("J126").Value = 3053.5
$ws.Range("K126").Value = 182256.18
$ws.Range("L126").Value = 9160.5
$ws.Range("M126").Value = -179786.18
$ws.Range("N126").Value = -14100.5

$ws = $wb.Sheets("LTW")
$ws.Range("H7").Value = 3175.3076
$ws.Range("I7").Value = 2276.25
$ws.Range("K7").Value = 2276.25
$ws.Range("M7").Value = -2164.25

$ws.Range("H9").Value = 370
$ws.Range("I9").Value = 55
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 55
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = 169
$ws.Range("N9").Value = -1448

$ws.Range("H16").Value = 1331.7307
$ws.Range("I16").Value = 1433.0588
$ws.Range("K16").Value = 1433.0588
$ws.Range("M16").Value = -1263.0588

$ws.Range("H40").Value = 3408.25
$ws.Range("I40").Value = 3650.5
$ws.Range("K40").Value = 3650.5
$ws.Range("M40").Value = -3514.5

$ws.Range("H122").Value = 3998.8572
$ws.Range("I122").Value = 3719.8
$ws.Range("K122").Value = 11159.4
$ws.Range("M122").Value = -8709.400000000001

$ws.Range("H126").Value = 3175.3076
$ws.Range("I126").Value = 2276.25
$ws.Range("K126").Value = 6828.75
$ws.Range("M126").Value = -4358.75

$ws = $wb.Sheets("WVR")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H33").Value = 15000
$ws.Range("I33").Value = 15000
$ws.Range("K33").Value = 15000
$ws.Range("M33").Value = -14750

$ws.Range("H36").Value = 15000
$ws.Range("I36").Value = 15000
$ws.Range("K36").Value = 15000
$ws.Range("M36").Value = -14750

$ws.Range("H122").Value = 103386.7
$ws.Range("I122").Value = 128301.5
$ws.Range("K122").Value = 384904.5
$ws.Range("M122").Value = -382454.5

$ws.Range("H126").Value = 10105.714
$ws.Range("I126").Value = 11590.5
$ws.Range("K126").Value = 34771.5
$ws.Range("M126").Value = -32301.5
